# Auto-generated edit script: regenerate validation-result data table
# (new simulation run with 5 angle groups instead of 4) and refresh the
# summary statistics block (H2:I6). Matches the "ErrorDetecter and
# PredictSA function added" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=S_Real  B=Angle  C=S_Pred  D=Angle_Pred  E=S_Error  F=Angle_Error
$data = @(
    @(2550.4,0,2550.161392067887,0.00638856355559393,-0.23860793211315467,0.00638856355559393),
    @(2550.4,-15.9,2550.5496096175857,-15.877600439564352,0.14960961758561098,0.02239956043564817),
    @(2550.4,-14.1,2550.6140396978662,-14.099999993248332,0.21403969786615562,0.000000006751667669391281),
    @(2550.4,-12.03,2549.6053711558693,-12.095788659944784,-0.7946288441307843,-0.06578865994478456),
    @(2550.4,-10.07,2550.639055106227,-10.1527801091414,0.239055106226715,-0.08278010914139955),
    @(2550.4,-8.07,2549.618903668694,-8.169329712330601,-0.7810963313058892,-0.09932971233060073),
    @(2550.4,-6.1,2550.17667080602,-6.166523917107877,-0.22332919398013473,-0.06652391710787775),
    @(2550.4,-4.1,2550.1778662328234,-4.148856783340104,-0.22213376717672872,-0.04885678334010457),
    @(2550.4,-2.08,2550.595855496092,-2.105615599294953,0.19585549609200825,-0.02561559929495294),
    @(2550.4,-0.08,2550.124466641725,-0.09207660091834034,-0.27553335827496994,-0.012076600918340341),
    @(2550.4,1.92,2550.414125408703,1.9563354161263553,0.014125408702966524,0.03633541612635538),
    @(2550.4,3.93,2550.255514480605,3.986911728778637,-0.14448551939494791,0.05691172877863693),
    @(2550.4,5.9,2550.103599253121,5.9959377855066665,-0.29640074687904416,0.0959377855066661),
    @(2550.4,7.93,2549.9143772912626,8.020461571085388,-0.4856227087375373,0.09046157108538821),
    @(2550.4,9.93,2549.907739152298,10.01473760949381,-0.492260847702255,0.0847376094938106),
    @(2550.4,11.95,2549.6076512331892,12.002012019321093,-0.7923487668108464,0.05201201932109356),
    @(2550.4,13.95,2550.1954008108905,13.949999999999996,-0.2045991891095582,-0.000000000000003552713678800501),
    @(2550.4,15.97,2551.3930021847286,15.905341518726285,0.9930021847285389,-0.06465848127371565),
    @(2650.1,0,2650.7969863636044,0.011775580175488633,0.6969863636045375,0.011775580175488633),
    @(2650.1,-15.9,2650.470405046899,-15.839827574461982,0.3704050468991227,0.060172425538018715),
    @(2650.1,-14.07,2650.146668536862,-14.07000000000002,0.04666853686194372,-0.000000000000019539925233402755),
    @(2650.1,-12.05,2650.4703638588944,-12.089278362936085,0.37036385889450685,-0.03927836293608422),
    @(2650.1,-10.07,2650.7510966173368,-10.134494271539328,0.6510966173368615,-0.06449427153932774),
    @(2650.1,-8.07,2651.1448838388123,-8.159281441898289,1.0448838388124386,-0.08928144189828835),
    @(2650.1,-6.1,2651.148826746839,-6.157607642478627,1.0488267468390404,-0.057607642478627596),
    @(2650.1,-4.1,2649.367181384447,-4.142943714026313,-0.7328186155527874,-0.04294371402631292),
    @(2650.1,-2.08,2650.3038113374046,-2.109403823212674,0.20381133740465884,-0.029403823212673874),
    @(2650.1,-0.08,2649.8851160040053,-0.09260958318582123,-0.21488399599456898,-0.012609583185821227),
    @(2650.1,1.92,2649.9222517575454,1.9494758309964166,-0.17774824245452692,0.029475830996416708),
    @(2650.1,3.93,2650.4273113187082,3.978179952265137,0.3273113187083254,0.048179952265136716),
    @(2650.1,5.9,2650.267916133782,5.978175704677262,0.16791613378200054,0.07817570467726132),
    @(2650.1,7.93,2649.9603972392683,7.996871202438898,-0.1396027607315773,0.06687120243889844),
    @(2650.1,9.93,2650.2925339419503,9.996537808923343,0.1925339419503871,0.06653780892334282),
    @(2650.1,11.95,2649.0593291611663,11.974027327246795,-1.0406708388336483,0.024027327246795238),
    @(2650.1,13.95,2649.5775204417755,13.936423646583327,-0.5224795582244042,-0.013576353416672404),
    @(2650.1,15.97,2651.7330811144134,15.870775618404785,1.6330811144134714,-0.0992243815952154),
    @(2749.8,0,2750.1266924735996,0.005248521283036148,0.32669247359945075,0.005248521283036148),
    @(2749.8,-15.9,2749.612978464957,-15.827120831179167,-0.1870215350431863,0.07287916882083323),
    @(2749.8,-14.1,2750.374625772989,-14.040863433933792,0.5746257729888384,0.05913656606620776),
    @(2749.8,-12.05,2748.993680779263,-12.03159475681788,-0.8063192207373504,0.01840524318212111),
    @(2749.8,-10.05,2750.081863311489,-10.100076709985428,0.2818633114889053,-0.05007670998542757),
    @(2749.8,-8.07,2749.983540461779,-8.135325663330956,0.18354046177864802,-0.06532566333095602),
    @(2749.8,-6.1,2750.3323073960705,-6.15716353218214,0.5323073960703368,-0.05716353218214021),
    @(2749.8,-4.1,2749.8723567699644,-4.134121922935293,0.07235676996424445,-0.03412192293529337),
    @(2749.8,-2.08,2750.504949161105,-2.101213873351505,0.7049491611046506,-0.02121387335150482),
    @(2749.8,-0.08,2750.2238749920307,-0.0834738392171428,0.4238749920305054,-0.003473839217142796),
    @(2749.8,1.92,2749.2751795898057,1.9481340868812422,-0.5248204101944793,0.028134086881242304),
    @(2749.8,3.93,2749.55534319062,3.971587433808035,-0.24465680938010337,0.041587433808034646),
    @(2749.8,5.9,2749.663167299467,5.965625051040167,-0.13683270053297747,0.06562505104016658),
    @(2749.8,7.93,2749.6535719574617,7.9839246474217145,-0.14642804253844588,0.053924647421714766),
    @(2749.8,9.93,2749.4415263723704,9.987445819202302,-0.35847362762979174,0.05744581920230196),
    @(2749.8,11.95,2749.4096650611405,11.949999999999996,-0.3903349388597235,-0.000000000000003552713678800501),
    @(2749.8,13.95,2750.079305047471,13.900990411348143,0.27930504747064333,-0.049009588651856006),
    @(2749.8,15.97,2750.7538588820753,15.840454279572628,0.9538588820751102,-0.1295457204273731),
    @(2849.5,0,2848.8390434001544,0.012403651808405029,-0.6609565998455764,0.012403651808405029),
    @(2849.5,-15.9,2849.809850008129,-15.790578766958992,0.3098500081291604,0.10942123304100804),
    @(2849.5,-14.1,2849.2623001431393,-14.022437498241038,-0.2376998568606723,0.07756250175896184),
    @(2849.5,-12.05,2848.6192479911906,-12.020854084794182,-0.8807520088093952,0.02914591520581844),
    @(2849.5,-10.05,2850.2966525972315,-10.088480624733293,0.7966525972315139,-0.0384806247332925),
    @(2849.5,-8.07,2848.6663673612584,-8.118359517700025,-0.8336326387416193,-0.04835951770002467),
    @(2849.5,-6.1,2849.1642076949543,-6.137786045407498,-0.335792305045743,-0.03778604540749875),
    @(2849.5,-4.1,2849.0032808360957,-4.125066156236186,-0.4967191639043449,-0.025066156236186465),
    @(2849.5,-2.08,2849.516174966876,-2.100344995338662,0.016174966875951213,-0.020344995338661853),
    @(2849.5,-0.08,2849.3108860659954,-0.09870601279288149,-0.1891139340045811,-0.01870601279288149),
    @(2849.5,1.92,2850.6000719304916,1.9362269908974676,1.1000719304915947,0.016226990897467708),
    @(2849.5,3.93,2850.6303774758876,3.961447055700325,1.1303774758875988,0.0314470557003248),
    @(2849.5,5.9,2849.0967386367124,5.962534728507994,-0.40326136328758366,0.0625347285079938),
    @(2849.5,7.93,2850.128029574524,7.966595148402985,0.6280295745241347,0.03659514840298517),
    @(2849.5,9.93,2849.3795356357955,9.94422755216903,-0.12046436420450846,0.014227552169030488),
    @(2849.5,11.95,2848.1059483446647,11.932431455255037,-1.3940516553352609,-0.017568544744962722),
    @(2849.5,13.95,2849.4360748626877,13.875677673147843,-0.06392513731225336,-0.07432232685215645),
    @(2849.5,15.97,2850.218980288692,15.816899981057274,0.7189802886919097,-0.15310001894272673),
    @(2949.2,0,2950.710447587979,0.014879995975111101,1.5104475879793426,0.014879995975111101),
    @(2949.2,-15.9,2948.6203887089787,-15.76559241100929,-0.5796112910211377,0.13440758899070993),
    @(2949.2,-14.1,2949.17921020076,-13.994067069332408,-0.020789799239992135,0.10593293066759202),
    @(2949.2,-12.05,2947.745583148298,-11.995935183795245,-1.454416851701808,0.05406481620475567),
    @(2949.2,-10.05,2949.6736251312977,-10.058192808356788,0.47362513129792205,-0.00819280835678704),
    @(2949.2,-8.07,2948.6141209947905,-8.109967637619615,-0.5858790052093354,-0.03996763761961475),
    @(2949.2,-6.1,2949.062625615876,-6.127740621917302,-0.1373743841236319,-0.0277406219173022),
    @(2949.2,-4.1,2948.284725546566,-4.111628755887804,-0.9152744534339945,-0.011628755887804765),
    @(2949.2,-2.08,2950.8164492952656,-2.0994156296541346,1.6164492952657383,-0.01941562965413457),
    @(2949.2,-0.08,2950.551106061373,-0.08000000000000813,1.3511060613732297,-0.000000000000008132383655379272),
    @(2949.2,1.92,2949.0575486294406,1.9336146689402647,-0.14245137055922896,0.013614668940264751),
    @(2949.2,3.93,2948.540075212161,3.964690263784027,-0.6599247878389178,0.03469026378402695),
    @(2949.2,5.9,2949.2252180391256,5.953082970618161,0.02521803912577525,0.05308297061816081),
    @(2949.2,7.93,2948.250249459354,7.958482788493167,-0.9497505406457094,0.0284827884931671),
    @(2949.2,9.93,2949.47799237456,9.926342888293687,0.2779923745601991,-0.0036571117063122216),
    @(2949.2,11.95,2948.4179006741388,11.91561618346365,-0.7820993258610542,-0.0343838165363497),
    @(2949.2,13.95,2948.1392235052276,13.852135181544416,-1.060776494772199,-0.0978648184555837),
    @(2949.2,15.97,2949.804914845433,15.797021999126999,0.6049148454330862,-0.17297800087300175)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Summary block (H2:I6): Mean_Error, Samples, Qualified, Qua_Rate, 95%_Confidence_Int
$ws.Cells.Item(2, 9).Value = 0.521463251624775
$ws.Cells.Item(3, 9).Value = 90
$ws.Cells.Item(4, 9).Value = 90
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(6, 9).Value = 1.3940516553352609

# Scroll the view right so column E is the left-most visible column
# (mirrors the author's window state after adding more rows of data),
# without disturbing the current H2:I6 selection.
try {
    $excel.ActiveWindow.ScrollRow = 1
    $excel.ActiveWindow.ScrollColumn = 5
} catch {
}
